$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.284.25'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '2.682.76'
$ws.Range('E3').Value = '  -2.85%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '554.16'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -3.85%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '157.92'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E8').Value = '  -3.44%  '
$ws.Range('E9').Value = '  -4.62%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.164'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -0.09%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.367'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -4.81%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '5.35'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -8.46%  '
$ws.Range('D13').Value = '3.161.25'
$ws.Range('E13').Value = '  -2.74%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '26.25'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -2.77%  '
$ws.Range('D15').Value = '63.154.55'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('E16').Value = '  -4.59%  '
$ws.Range('D17').Value = '2.690.84'
$ws.Range('E17').Value = '  -2.77%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '12.01'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.90%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '4.56'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -5.86%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '342.26'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -4.88%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.30'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -5.51%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.994'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -0.48%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '0.503'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -5.01%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '63.36'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -2.53%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '0.168'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -1.04%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -0.27%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '8.02'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -5.90%  '
$ws.Range('D28').Value = '0.0₃0851'
$ws.Range('E28').Value = '  -6.08%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.92'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -2.01%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.32'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +4.60%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '7.00'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -4.92%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '165.24'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -3.07%  '
$ws.Range('E33').Value = '  +0.01%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '19.50'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -3.47%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '4.72'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -3.92%  '
$ws.Range('E36').Value = '  -3.04%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '1.77'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -2.34%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '338.76'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -2.79%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.935'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -6.68%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '6.08'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -2.80%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '3.93'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -5.53%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '37.99'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '20.22'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -5.75%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '20.69'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -5.17%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.618'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -1.89%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.0559'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -4.83%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +0.03%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '11.05'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -0.12%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '129.58'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -5.29%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0968'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -4.23%  '
$ws.Range('D51').Value = '2.093.44'
$ws.Range('E51').Value = '  -1.39%  '
